$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '27.836.68'
$ws.Range("E2").Value = '  +2.66%  '
$ws.Range("D3").Value = '1.768.59'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  -0.48%  '
Set-TextValue "D5" '335.03'
$ws.Range("E5").Value = '  -0.49%  '
Set-TextValue "D6" '0.9967'
$ws.Range("E6").Value = '  -0.60%  '
Set-TextValue "D7" '0.3812'
$ws.Range("E7").Value = '  -0.11%  '
Set-TextValue "D8" '0.3420'
$ws.Range("E8").Value = '  +0.35%  '
Set-TextValue "D9" '47.95'
$ws.Range("E9").Value = '  -0.13%  '
Set-TextValue "D10" '1.139'
$ws.Range("E10").Value = '  -3.86%  '
Set-TextValue "D11" '0.07400'
$ws.Range("E11").Value = '  -0.30%  '
Set-TextValue "D12" '0.9958'
$ws.Range("E12").Value = '  -0.77%  '
Set-TextValue "D13" '22.58'
$ws.Range("E13").Value = '  +4.58%  '
Set-TextValue "D14" '6.349'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").Value = '1.768.47'
$ws.Range("E15").Value = '  -0.67%  '
Set-TextValue "D16" '7.086'
$ws.Range("E16").Value = '  +0.60%  '
Set-TextValue "D17" '0.00001076'
$ws.Range("E17").Value = '  -0.57%  '
Set-TextValue "D18" '0.06668'
$ws.Range("E18").Value = '  +0.44%  '
Set-TextValue "D19" '82.05'
$ws.Range("E19").Value = '  -1.44%  '
Set-TextValue "D20" '0.9974'
$ws.Range("E20").Value = '  -0.53%  '
Set-TextValue "D21" '17.33'
$ws.Range("E21").Value = '  +0.37%  '
Set-TextValue "D22" '6.411'
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").Value = '27.854.79'
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("E24").Value = '  -1.14%  '
Set-TextValue "D25" '2.382'
$ws.Range("E25").Value = '  +0.37%  '
Set-TextValue "D26" '1.438'
$ws.Range("E26").Value = '  -1.16%  '
Set-TextValue "D27" '20.72'
$ws.Range("E27").Value = '  -1.50%  '
Set-TextValue "D28" '2.414'
$ws.Range("E28").Value = '  -3.47%  '
Set-TextValue "D29" '153.59'
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("D30").Value = '1.969.78'
$ws.Range("E30").Value = '  -0.48%  '
Set-TextValue "D31" '134.19'
$ws.Range("E31").Value = '  +0.28%  '
Set-TextValue "D32" '6.142'
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("E34").Value = '  +1.40%  '
Set-TextValue "D35" '12.75'
$ws.Range("E35").Value = '  -2.24%  '
Set-TextValue "D36" '0.02430'
Set-TextValue "D37" '0.6851'
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("E38").Value = '  -1.20%  '
Set-TextValue "D39" '0.06320'
$ws.Range("E39").Value = '  +0.87%  '
Set-TextValue "D40" '0.2184'
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  -6.43%  '
Set-TextValue "D42" '1.234'
$ws.Range("E42").Value = '  +0.58%  '
Set-TextValue "D43" '8.234'
$ws.Range("E43").Value = '  -3.49%  '
Set-TextValue "D44" '14.12'
$ws.Range("E44").Value = '  -0.80%  '
Set-TextValue "D45" '0.9966'
$ws.Range("E45").Value = '  -0.60%  '
Set-TextValue "D46" '0.6267'
$ws.Range("E46").Value = '  -2.11%  '
Set-TextValue "D47" '3.838'
$ws.Range("E47").Value = '  -0.50%  '
Set-TextValue "D48" '131.38'
$ws.Range("E48").Value = '  +0.18%  '
Set-TextValue "D49" '2.080'
$ws.Range("E49").Value = '  -1.64%  '
Set-TextValue "D50" '0.07371'
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue "D51" '1.144'
$ws.Range("E51").Value = '  +2.71%  '
